$wb = $excel.ActiveWorkbook

# --- Worksheet references ---
$wsCurves = $wb.Worksheets.Item("curves")
$wsGroups = $wb.Worksheets.Item("groups")
$wsG2c    = $wb.Worksheets.Item("g2c")
$wsC2t    = $wb.Worksheets.Item("c2t")

# --- Text updates (hyphens -> spaces, and ADHD / Aspergers -> ADHD/Aspergers) ---
# Order matters: the first time each new distinct string is written it becomes a
# freshly appended shared-string entry, so keep this ordering to match the
# canonical shared string table ordering produced by the edit.

# "Self harm" (was "Self-harm")
$wsG2c.Range("B6").Value = "Self harm"
$wsG2c.Range("B26").Value = "Self harm"
$wsG2c.Range("B57").Value = "Self harm"
$wsC2t.Range("A44").Value = "Self harm"
$wsC2t.Range("A45").Value = "Self harm"
$wsC2t.Range("A46").Value = "Self harm"
$wsC2t.Range("A47").Value = "Self harm"
$wsC2t.Range("A48").Value = "Self harm"
$wsC2t.Range("A49").Value = "Self harm"
$wsC2t.Range("A50").Value = "Self harm"

# "Family of COVID deceased" (was "Family of COVID-deceased")
$wsGroups.Range("A5").Value = "Family of COVID deceased"
$wsG2c.Range("A16").Value = "Family of COVID deceased"
$wsG2c.Range("A17").Value = "Family of COVID deceased"

# "Pre existing CMH illness" (was "Pre-existing CMH illness")
$wsGroups.Range("A13").Value = "Pre existing CMH illness"
$wsG2c.Range("A48").Value = "Pre existing CMH illness"

# "Pre existing LTC" (was "Pre-existing LTC")
$wsGroups.Range("A14").Value = "Pre existing LTC"
$wsG2c.Range("A49").Value = "Pre existing LTC"
$wsG2c.Range("A50").Value = "Pre existing LTC"

# "Pre existing SMI" (was "Pre-existing SMI")
$wsGroups.Range("A15").Value = "Pre existing SMI"
$wsG2c.Range("A51").Value = "Pre existing SMI"

# "Shallow mid term" (was "Shallow mid-term")
$wsCurves.Range("D1").Value = "Shallow mid term"

# "Neurological symptom disorder (ADHD/Aspergers)" (was "... (ADHD / Aspergers)")
$wsG2c.Range("B4").Value = "Neurological symptom disorder (ADHD/Aspergers)"
$wsC2t.Range("A25").Value = "Neurological symptom disorder (ADHD/Aspergers)"
$wsC2t.Range("A26").Value = "Neurological symptom disorder (ADHD/Aspergers)"

# --- View state changes ---
# Previously "curves" was the selected tab with default A1 selection; now "c2t" is
# the selected tab with A27 selected, and "g2c" is scrolled down with B5 selected.

$wsG2c.Activate()
$wsG2c.Range("B5").Select()

$wsC2t.Activate()
$wsC2t.Range("A27").Select()
